$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing (stale) hyperlink pointing at the old F2 / Maryem record
$ws.Hyperlinks.Delete()

# Remove the table/list object (registro_docente) - converts back to a normal range
$ws.ListObjects.Item(1).Unlist()

# Delete rows 2 (Maryem Ruiz, test data) and 3 ("ss" placeholder row),
# leaving what was row 4 (Sebastian Palacio) as the new row 2
$ws.Rows("2:3").Delete()

# Header row should not retain the centered table-header style anymore
$ws.Range("A1:J1").Style = "Normal"

# Fix the typo in Sebastian's e-mail address (comma -> dots)
$ws.Range("F2").Value = "sebastian_palacio23231@elpoli.edu.co"

# Re-create the hyperlink on the corrected e-mail cell
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:sebastian_palacio23231@elpoli.edu.co") | Out-Null

# Reset selection to the natural top-left cell
$ws.Range("A1").Select() | Out-Null
